$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0,
    0.1288322815721658,
    0.2657586058018112,
    0.4193815332872734,
    0.5993526087583387,
    -0.4434531441505186,
    -0.4434531441505186,
    -0.6013775875253951,
    -0.4743415942963839,
    -0.4743415942963839,
    0.001834528240272033,
    0.001834528240272033,
    0.2473474482405971,
    0.1063814174549735,
    0.1063814174549735,
    -0.02285419592424004,
    -0.02285419592424004,
    0.08759328731672891,
    0.188609191195769,
    0.188609191195769,
    -0.02285419592424004,
    -0.02285419592424004,
    0.08759328731672891,
    0.188609191195769,
    0.188609191195769,
    2.71,
    2.779108819735061,
    2.990777123454038,
    3.358303226289058,
    3.904777396593537,
    1.644156523126431,
    1.612743423246858,
    1.829189117262751,
    1.828540277365283,
    1.995597596835747,
    1.322189866823892,
    1.225977001497085,
    0.9792252696219189,
    1.186055825661855,
    1.282268690988662,
    1.544110198463599,
    1.377052878993135,
    1.204843955800101,
    1.112526171453693,
    1.143939271333267,
    1.792507548601998,
    1.544110198463599,
    1.301056821126909,
    1.143939271333267,
    1.112526171453693
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
